$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells in column D whose new price text would otherwise be auto-detected as a
# number by Excel (losing trailing zeros / exact formatting). Mark them as Text
# first so the literal string is preserved, one cell at a time (a multi-area
# Range only keeps the NumberFormat override on its first area).
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"

# Updated Price (D) / Volume(1h) (E) figures
$ws.Range("D2").Value = "27.118.33"
$ws.Range("D3").Value = "1.866.11"
$ws.Range("E3").Value = "  -2.17%  "
$ws.Range("D4").Value = "1.000"
$ws.Range("E4").Value = "  -0.09%  "
$ws.Range("D5").Value = "306.33"
$ws.Range("E5").Value = "  -2.07%  "
$ws.Range("E6").Value = "  +0.02%  "
$ws.Range("D7").Value = "0.5147"
$ws.Range("E7").Value = "  -1.54%  "
$ws.Range("D8").Value = "0.3764"
$ws.Range("E8").Value = "  -0.53%  "
$ws.Range("D9").Value = "0.07156"
$ws.Range("E9").Value = "  -1.20%  "
$ws.Range("E10").Value = "  -1.80%  "
$ws.Range("D11").Value = "20.74"
$ws.Range("E11").Value = "  -2.83%  "
$ws.Range("D12").Value = "0.07564"
$ws.Range("E12").Value = "  -1.26%  "
$ws.Range("D13").Value = "1.842.91"
$ws.Range("E13").Value = "  -3.54%  "
$ws.Range("D14").Value = "5.310"
$ws.Range("E14").Value = "  -2.69%  "
$ws.Range("D15").Value = "89.72"
$ws.Range("E15").Value = "  -2.54%  "
$ws.Range("E16").Value = "  -0.11%  "
$ws.Range("D17").Value = "0.000008473"
$ws.Range("E17").Value = "  -2.79%  "
$ws.Range("D18").Value = "14.05"
$ws.Range("E18").Value = "  -3.35%  "
$ws.Range("E19").Value = "  -0.08%  "
$ws.Range("D20").Value = "27.146.69"
$ws.Range("E21").Value = "  -2.74%  "
$ws.Range("D22").Value = "2.085.43"
$ws.Range("E22").Value = "  -3.24%  "
$ws.Range("E23").Value = "  -3.44%  "
$ws.Range("D24").Value = "6.446"
$ws.Range("E24").Value = "  -2.93%  "
$ws.Range("D25").Value = "1.838"
$ws.Range("E25").Value = "  -1.80%  "
$ws.Range("D26").Value = "146.20"
$ws.Range("E26").Value = "  -4.97%  "
$ws.Range("D27").Value = "17.96"
$ws.Range("E27").Value = "  -2.23%  "
$ws.Range("D28").Value = "2.093"
$ws.Range("E28").Value = "  -3.45%  "
$ws.Range("D29").Value = "112.80"
$ws.Range("E29").Value = "  -1.86%  "
$ws.Range("D30").Value = "4.666"
$ws.Range("E30").Value = "  -3.95%  "
$ws.Range("D31").Value = "4.658"
$ws.Range("E31").Value = "  -4.21%  "
$ws.Range("E32").Value = "  +0.70%  "
$ws.Range("E33").Value = "  -3.28%  "
$ws.Range("E34").Value = "  -3.47%  "
$ws.Range("D36").Value = "0.7250"
$ws.Range("E36").Value = "  -7.27%  "
$ws.Range("D37").Value = "0.02037"
$ws.Range("E37").Value = "  -2.86%  "
$ws.Range("D38").Value = "3.086"
$ws.Range("E38").Value = "  +0.29%  "
$ws.Range("D39").Value = "2.497"
$ws.Range("E39").Value = "  -4.27%  "
$ws.Range("D40").Value = "1.076"
$ws.Range("E40").Value = "  -1.48%  "
$ws.Range("D41").Value = "0.5284"
$ws.Range("E41").Value = "  -5.75%  "
$ws.Range("D42").Value = "6.466"
$ws.Range("E42").Value = "  -3.80%  "
$ws.Range("D43").Value = "115.77"
$ws.Range("E43").Value = "  +0.09%  "
$ws.Range("D44").Value = "8.278"
$ws.Range("E44").Value = "  -3.48%  "
$ws.Range("E45").Value = "  -3.50%  "
$ws.Range("D46").Value = "1.000"
$ws.Range("E46").Value = "  +0.04%  "
$ws.Range("D47").Value = "0.4621"
$ws.Range("E47").Value = "  -4.10%  "
$ws.Range("D48").Value = "9.940"
$ws.Range("E48").Value = "  -5.51%  "
$ws.Range("D49").Value = "1.565"
$ws.Range("E49").Value = "  -3.50%  "
$ws.Range("D50").Value = "36.53"
$ws.Range("E50").Value = "  -1.44%  "
$ws.Range("D51").Value = "63.45"
$ws.Range("E51").Value = "  -5.24%  "

# Remove the temporary Text number-format so the cells end up styled exactly
# like the rest of the column again (no explicit style index).
$ws.Range("D4").ClearFormats()
$ws.Range("D5").ClearFormats()
$ws.Range("D7").ClearFormats()
$ws.Range("D8").ClearFormats()
$ws.Range("D9").ClearFormats()
$ws.Range("D11").ClearFormats()
$ws.Range("D12").ClearFormats()
$ws.Range("D14").ClearFormats()
$ws.Range("D15").ClearFormats()
$ws.Range("D17").ClearFormats()
$ws.Range("D18").ClearFormats()
$ws.Range("D24").ClearFormats()
$ws.Range("D25").ClearFormats()
$ws.Range("D26").ClearFormats()
$ws.Range("D27").ClearFormats()
$ws.Range("D28").ClearFormats()
$ws.Range("D29").ClearFormats()
$ws.Range("D30").ClearFormats()
$ws.Range("D31").ClearFormats()
$ws.Range("D36").ClearFormats()
$ws.Range("D37").ClearFormats()
$ws.Range("D38").ClearFormats()
$ws.Range("D39").ClearFormats()
$ws.Range("D40").ClearFormats()
$ws.Range("D41").ClearFormats()
$ws.Range("D42").ClearFormats()
$ws.Range("D43").ClearFormats()
$ws.Range("D44").ClearFormats()
$ws.Range("D46").ClearFormats()
$ws.Range("D47").ClearFormats()
$ws.Range("D48").ClearFormats()
$ws.Range("D49").ClearFormats()
$ws.Range("D50").ClearFormats()
$ws.Range("D51").ClearFormats()
